$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, [string]$text) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
    $cellRange.ClearFormats()
}

Set-TextValue $ws.Range("D2") "295.25"
Set-TextValue $ws.Range("E2") "0.46%"
Set-TextValue $ws.Range("G2") "20"
Set-TextValue $ws.Range("D3") "31.54"
Set-TextValue $ws.Range("E3") "1.66%"
Set-TextValue $ws.Range("G3") "20"
Set-TextValue $ws.Range("D4") "4.966"
Set-TextValue $ws.Range("E4") "0.37%"
Set-TextValue $ws.Range("G4") "20"
Set-TextValue $ws.Range("D5") "0.07619"
Set-TextValue $ws.Range("E5") "3.82%"
Set-TextValue $ws.Range("G5") "20"
Set-TextValue $ws.Range("D6") "2.256"
Set-TextValue $ws.Range("E6") "-5.63%"
Set-TextValue $ws.Range("G6") "20"
Set-TextValue $ws.Range("D7") "7.839"
Set-TextValue $ws.Range("E7") "1.57%"
Set-TextValue $ws.Range("G7") "20"
Set-TextValue $ws.Range("B8") "MXToken"
Set-TextValue $ws.Range("C8") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D8") "0.9224"
Set-TextValue $ws.Range("E8") "2.30%"
Set-TextValue $ws.Range("G8") "20"
Set-TextValue $ws.Range("B9") "LiechtensteinCryptoassetsExchange"
Set-TextValue $ws.Range("C9") "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D9") "0.09616"
Set-TextValue $ws.Range("E9") "21.01%"
Set-TextValue $ws.Range("G9") "20"
Set-TextValue $ws.Range("B10") "WazirX"
Set-TextValue $ws.Range("C10") "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D10") "0.1738"
Set-TextValue $ws.Range("E10") "3.49%"
Set-TextValue $ws.Range("G10") "20"
Set-TextValue $ws.Range("B11") "MandalaExchangeToken"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D11") "0.08389"
Set-TextValue $ws.Range("E11") "3.01%"
Set-TextValue $ws.Range("G11") "20"
Set-TextValue $ws.Range("B12") "BitrueCoin"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D12") "0.03251"
Set-TextValue $ws.Range("E12") "4.81%"
Set-TextValue $ws.Range("G12") "20"
Set-TextValue $ws.Range("B13") "BitMartToken"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D13") "0.09977"
Set-TextValue $ws.Range("E13") "-1.05%"
Set-TextValue $ws.Range("G13") "20"
Set-TextValue $ws.Range("B14") "BitForexToken"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D14") "0.001483"
Set-TextValue $ws.Range("E14") "-1.53%"
Set-TextValue $ws.Range("G14") "20"
Set-TextValue $ws.Range("B15") "TigerCash"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D15") "0.005687"
Set-TextValue $ws.Range("E15") "-2.09%"
Set-TextValue $ws.Range("G15") "20"
Set-TextValue $ws.Range("B16") "LEO"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D16") "3.484"
Set-TextValue $ws.Range("E16") "0.07%"
Set-TextValue $ws.Range("G16") "20"
Set-TextValue $ws.Range("B17") "GateToken"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D17") "3.773"
Set-TextValue $ws.Range("E17") "0.88%"
Set-TextValue $ws.Range("G17") "20"
Set-TextValue $ws.Range("D18") "2.142"
Set-TextValue $ws.Range("E18") "3.19%"
Set-TextValue $ws.Range("G18") "20"
Set-TextValue $ws.Range("E19") "0.57%"
Set-TextValue $ws.Range("G19") "20"
Set-TextValue $ws.Range("D20") "0.1318"
Set-TextValue $ws.Range("E20") "1.25%"
Set-TextValue $ws.Range("G20") "20"
Set-TextValue $ws.Range("D21") "4.072"
Set-TextValue $ws.Range("E21") "2.32%"
Set-TextValue $ws.Range("G21") "20"
Set-TextValue $ws.Range("D22") "0.2272"
Set-TextValue $ws.Range("E22") "8.35%"
Set-TextValue $ws.Range("G22") "20"
Set-TextValue $ws.Range("D23") "0.04512"
Set-TextValue $ws.Range("E23") "-0.70%"
Set-TextValue $ws.Range("G23") "20"
Set-TextValue $ws.Range("D24") "0.001208"
Set-TextValue $ws.Range("E24") "-0.14%"
Set-TextValue $ws.Range("G24") "20"
Set-TextValue $ws.Range("D25") "0.004329"
Set-TextValue $ws.Range("E25") "-6.84%"
Set-TextValue $ws.Range("G25") "20"
Set-TextValue $ws.Range("D26") "0.0001292"
Set-TextValue $ws.Range("E26") "-0.40%"
Set-TextValue $ws.Range("G26") "20"
Set-TextValue $ws.Range("D27") "0.0003360"
Set-TextValue $ws.Range("E27") "-0.86%"
Set-TextValue $ws.Range("G27") "20"
Set-TextValue $ws.Range("G28") "20"
Set-TextValue $ws.Range("G29") "20"
Set-TextValue $ws.Range("G30") "20"
Set-TextValue $ws.Range("G31") "20"
Set-TextValue $ws.Range("G32") "20"
Set-TextValue $ws.Range("G33") "20"
Set-TextValue $ws.Range("G34") "20"
Set-TextValue $ws.Range("G35") "20"
Set-TextValue $ws.Range("G36") "20"
Set-TextValue $ws.Range("G37") "20"
Set-TextValue $ws.Range("G38") "20"
Set-TextValue $ws.Range("D39") "0.01673"
Set-TextValue $ws.Range("E39") "4.01%"
Set-TextValue $ws.Range("G39") "20"
Set-TextValue $ws.Range("D40") "0.04613"
Set-TextValue $ws.Range("E40") "3.64%"
Set-TextValue $ws.Range("G40") "20"
Set-TextValue $ws.Range("D41") "0.007473"
Set-TextValue $ws.Range("E41") "1.57%"
Set-TextValue $ws.Range("G41") "20"
Set-TextValue $ws.Range("D42") "0.009731"
Set-TextValue $ws.Range("E42") "13.10%"
Set-TextValue $ws.Range("G42") "20"
Set-TextValue $ws.Range("D43") "0.1373"
Set-TextValue $ws.Range("E43") "3.50%"
Set-TextValue $ws.Range("G43") "20"
Set-TextValue $ws.Range("D44") "0.002143"
Set-TextValue $ws.Range("E44") "7.29%"
Set-TextValue $ws.Range("G44") "20"
Set-TextValue $ws.Range("D45") "0.009427"
Set-TextValue $ws.Range("E45") "-0.63%"
Set-TextValue $ws.Range("G45") "20"
Set-TextValue $ws.Range("D46") "0.00006039"
Set-TextValue $ws.Range("E46") "1.92%"
Set-TextValue $ws.Range("G46") "20"
Set-TextValue $ws.Range("D47") "0.00000000742"
Set-TextValue $ws.Range("E47") "-0.80%"
Set-TextValue $ws.Range("G47") "20"
Set-TextValue $ws.Range("D48") "2.551"
Set-TextValue $ws.Range("E48") "13.83%"
Set-TextValue $ws.Range("G48") "20"
Set-TextValue $ws.Range("D49") "0.001980"
Set-TextValue $ws.Range("E49") "-31.59%"
Set-TextValue $ws.Range("G49") "20"
Set-TextValue $ws.Range("D50") "0.00002079"
Set-TextValue $ws.Range("E50") "-0.80%"
Set-TextValue $ws.Range("G50") "20"
Set-TextValue $ws.Range("D51") "0.0001980"
Set-TextValue $ws.Range("E51") "-0.80%"
Set-TextValue $ws.Range("G51") "20"
